$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select full columns A:B to match the recorded selection state
[void]$ws.Range("A1:B1048576").Select()

# Column B (row 2 .. row 152) target values, in order
$bValues = @(0.984375,0.84375,0.53125,0.484375,0.375,0.328125,0.40625,0.3125,0.390625,0.328125,0.28125,0.375,0.296875,0.3125,0.28125,0.296875,0.34375,0.25,0.359375,0.265625,0.28125,0.28125,0.21875,0.28125,0.265625,0.15625,0.21875,0.203125,0.328125,0.375,0.21875,0.515625,0.46875,0.359375,0.390625,0.34375,0.40625,0.375,0.328125,0.359375,0.28125,0.234375,0.21875,0.234375,0.3125,0.265625,0.1875,0.265625,0.265625,0.25,0.203125,0.25,0.203125,0.234375,0.25,0.203125,0.28125,0.21875,0.25,0.3125,0.1875,0.203125,0.234375,0.28125,0.203125,0.15625,0.203125,0.234375,0.21875,0.171875,0.203125,0.171875,0.203125,0.15625,0.171875,0.1875,0.203125,0.1875,0.203125,0.140625,0.171875,0.171875,0.203125,0.15625,0.1875,0.171875,0.140625,0.203125,0.234375,0.140625,0.171875,0.1875,0.1875,0.171875,0.125,0.1875,0.1875,0.1875,0.203125,0.203125,0.203125,0.1875,0.046875,0.171875,0.171875,0.125,0.15625,0.203125,0.09375,0.21875,0.140625,0.171875,0.125,0.109375,0.140625,0.1875,0.1875,0.265625,0.15625,0.1875,0.09375,0.15625,0.140625,0.1875,0.078125,0.21875,0.1875,0.125,0.109375,0.21875,0.171875,0.140625,0.1875,0.078125,0.234375,0.203125,0.09375,0.125,0.15625,0.109375,0.109375,0.140625,0.109375,0.1875,0.21875,0.15625,0.25,0.15625,0.09375,0.140625,0.1090909090909091)

$numRows = $bValues.Count
$startRow = 2

$data = New-Object 'object[,]' $numRows,2

for ($i = 0; $i -lt $numRows; $i++) {
    $row = $startRow + $i
    if ($row -le 101) {
        $data[$i,0] = $row - 2
    } else {
        $data[$i,0] = "<__main__.DisplayOutputs object at 0x7f78b1067c70>"
    }
    $data[$i,1] = $bValues[$i]
}

$endRow = $startRow + $numRows - 1
$rng = $ws.Range("A$startRow" + ":B$endRow")
$rng.Value2 = $data

Write-Host "Wrote rows $startRow..$endRow"
